# Updated cryptos list data (Price + Volume(1h) columns) on the "Sheet1" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are stored as text in the source data (e.g. "64.951.04"
# uses dots as thousands separators), so force text formatting before assigning
# the value -- otherwise Excel would auto-convert plain-decimal-looking strings
# (like "0.999") into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.951.04"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.153.31"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.81"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.57"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.150.49"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("E13").Value = "  +13.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.02"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.663.24"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.972.18"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.154.22"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").Value = "  +4.21%  "
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.25"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.80"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.23"
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.71"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.18"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.80"
$ws.Range("E28").Value = "  +6.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("E30").Value = "  +7.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.58"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.16"
$ws.Range("E34").Value = "  +7.33%  "
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.76"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0897"
$ws.Range("E37").Value = "  +9.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "463.62"
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  +7.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.64"
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.048.96"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  +7.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.38"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0583"
$ws.Range("E47").Value = "  +12.10%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.54"
$ws.Range("E51").Value = "  +1.38%  "
